# Adding AM and PM suffixes to the time values in column C of the
# "google calender formated sheet" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most times simply get " AM" appended; the two evening entries (18:45 and
# 22:00) get converted to their 12-hour "H:MM PM" equivalents.
$ws.Range("C2").Value  = "05:45 AM"
$ws.Range("C3").Value  = "06:10 AM"
$ws.Range("C4").Value  = "06:15 AM"
$ws.Range("C5").Value  = "06:20 AM"
$ws.Range("C6").Value  = "06:22 AM"
$ws.Range("C7").Value  = "06:25 AM"
$ws.Range("C8").Value  = "06:35 AM"
$ws.Range("C9").Value  = "06:40 AM"
$ws.Range("C10").Value = "06:40 AM"
$ws.Range("C11").Value = "07:50 AM"
$ws.Range("C12").Value = "08:25 AM"
$ws.Range("C13").Value = "08:35 AM"
$ws.Range("C14").Value = "08:40 AM"
$ws.Range("C15").Value = "05:30 AM"
$ws.Range("C16").Value = "06:15 AM"
$ws.Range("C17").Value = "07:30 AM"
$ws.Range("C18").Value = "12:30 AM"
$ws.Range("C19").Value = "06:00 AM"
$ws.Range("C20").Value = "08:15 AM"
$ws.Range("C21").Value = "6:45 PM"
$ws.Range("C24").Value = "10:00 PM"
